$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the data rows 16-19:
# Before: row16=(1100399613, JEAN CARLOS..., 1812, 1356800)
#         row17=(1100399613, JEAN CARLOS..., 1811, 1356800)
#         row18=(1100399613, JEAN CARLOS..., 1810, 1356800)
#         row19=(1049534523, JUNIOR JOSE..., 1809, 1400000)
# After:  row16=(1049534523, JUNIOR JOSE..., 1809, 1400000)
#         row17=(1100399613, JEAN CARLOS..., 1810, 1356800)
#         row18=(1100399613, JEAN CARLOS..., 1811, 1356800)
#         row19=(1100399613, JEAN CARLOS..., 1812, 1356800)

$ws.Range("C16").Value = "1049534523"
$ws.Range("D16").Value = "JUNIOR JOSE SALAS ARRIETA"
$ws.Range("E16").Value = "1809"
$ws.Range("G16").Value = 1400000

$ws.Range("C17").Value = "1100399613"
$ws.Range("D17").Value = "JEAN CARLOS ROBLES ATENCIA"
$ws.Range("E17").Value = "1810"
$ws.Range("G17").Value = 1356800

$ws.Range("C18").Value = "1100399613"
$ws.Range("D18").Value = "JEAN CARLOS ROBLES ATENCIA"
$ws.Range("E18").Value = "1811"
$ws.Range("G18").Value = 1356800

$ws.Range("C19").Value = "1100399613"
$ws.Range("D19").Value = "JEAN CARLOS ROBLES ATENCIA"
$ws.Range("E19").Value = "1812"
$ws.Range("G19").Value = 1356800
